$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format D2:E51 as Text before writing, to avoid numeric auto-conversion
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.988.64'
$ws.Range('E2').Value = '  +1.24%  '

$ws.Range('D3').Value = '1.641.35'
$ws.Range('E3').Value = '  +0.43%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '213.58'
$ws.Range('E5').Value = '  +0.65%  '

$ws.Range('E6').Value = '  +0.11%  '

$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').Value = '23.62'
$ws.Range('E8').Value = '  +1.37%  '

$ws.Range('E9').Value = '  -1.17%  '

$ws.Range('E10').Value = '  +0.44%  '

$ws.Range('D11').Value = '0.0879'
$ws.Range('E11').Value = '  +2.50%  '

$ws.Range('D12').Value = '1.874.48'
$ws.Range('E12').Value = '  +0.47%  '

$ws.Range('D13').Value = '1.643.59'
$ws.Range('E13').Value = '  +0.69%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '4.10'
$ws.Range('E14').Value = '  +1.33%  '

$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '0.575'
$ws.Range('E15').Value = '  +4.07%  '

$ws.Range('D16').Value = '65.93'
$ws.Range('E16').Value = '  +1.15%  '

$ws.Range('D17').Value = '27.983.68'
$ws.Range('E17').Value = '  +1.27%  '

$ws.Range('D18').Value = '232.64'
$ws.Range('E18').Value = '  +1.07%  '

$ws.Range('E19').Value = '  +0.50%  '

$ws.Range('D20').Value = '7.61'
$ws.Range('E20').Value = '  +0.51%  '

$ws.Range('E21').Value = '  +0.11%  '

$ws.Range('E22').Value = '  +1.24%  '

$ws.Range('E23').Value = '  +0.01%  '

$ws.Range('E24').Value = '  -0.78%  '

$ws.Range('D25').Value = '151.64'
$ws.Range('E25').Value = '  +1.82%  '

$ws.Range('D26').Value = '6.98'
$ws.Range('E26').Value = '  +1.31%  '

$ws.Range('E27').Value = '  +1.48%  '

$ws.Range('E28').Value = '  +0.01%  '

$ws.Range('E29').Value = '  +0.10%  '

$ws.Range('E31').Value = '  +0.20%  '

$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.416.21'
$ws.Range('E33').Value = '  -4.23%  '

$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '3.11'
$ws.Range('E34').Value = '  +0.68%  '

$ws.Range('E35').Value = '  +2.34%  '

$ws.Range('E36').Value = '  +0.93%  '

$ws.Range('D37').Value = '0.887'
$ws.Range('E37').Value = '  +1.03%  '

$ws.Range('D38').Value = '0.0168'

$ws.Range('E39').Value = '  +0.03%  '

$ws.Range('D40').Value = '0.916'
$ws.Range('E40').Value = '  -4.40%  '

$ws.Range('E41').Value = '  -0.03%  '

$ws.Range('E42').Value = '  +0.04%  '

$ws.Range('D43').Value = '1.87'
$ws.Range('E43').Value = '  +7.28%  '

$ws.Range('D44').Value = '66.29'
$ws.Range('E44').Value = '  -2.15%  '

$ws.Range('D45').Value = '5.48'
$ws.Range('E45').Value = '  +2.93%  '

$ws.Range('E46').Value = '  +0.42%  '

$ws.Range('D47').Value = '1.783.21'
$ws.Range('E47').Value = '  +0.52%  '

$ws.Range('D48').Value = '88.05'
$ws.Range('E48').Value = '  +0.64%  '

$ws.Range('E49').Value = '  +1.43%  '

$ws.Range('E50').Value = '  +0.40%  '

$ws.Range('E51').Value = '  -1.07%  '

# Restore default style (removes the temporary Text number format marker)
$dataRange.Style = "Normal"
